# "Generate Report for handback"
#
# The source .md file and its .xlf handoff have now actually been handed
# back, so:
#   - Status goes from "Not yet handed off" to "Handed back" (shown on the
#     Overview sheet as well as on each language sheet).
#   - Each language sheet's row for that file gets its "Latest Target
#     File" / "Latest Handback File" columns (E/F) filled in, pointing at
#     the same handed-off file / translated package as columns A/C.
#   - The "Latest Handback DateTime" column (G) is stamped with the real
#     handback time instead of the epoch placeholder.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZh       = $wb.Worksheets.Item("zh-cn")
$wsDe       = $wb.Worksheets.Item("de-de")

# --- Status: "Not yet handed off" -> "Handed back" everywhere it shows up.
$wsOverview.Range("B2").Value = "Handed back"
$wsOverview.Range("C2").Value = "Handed back"
$wsZh.Range("B2").Value = "Handed back"
$wsDe.Range("B2").Value = "Handed back"

# --- zh-cn: record the handback (target file / handback file / datetime).
$zhMdLink = $null
$zhXlfLink = $null
foreach ($h in $wsZh.Hyperlinks) {
    $addr = $h.Range.Address()
    if ($addr -eq '$A$2') { $zhMdLink = $h.Address }
    if ($addr -eq '$C$2') { $zhXlfLink = $h.Address }
}

$wsZh.Hyperlinks.Add($wsZh.Range("E2"), $zhMdLink, $null, $null, "1ddd3d0e-7b9e-45da-b7a2-ededd37b75c7.md")
$wsZh.Hyperlinks.Add($wsZh.Range("F2"), $zhXlfLink, $null, $null, "1ddd3d0e-7b9e-45da-b7a2-ededd37b75c7.07b5d910097a9d19b5106c498aa4f82131c0eb33.zh-cn.xlf")
$wsZh.Range("G2").Value = "2016-01-08 09:20:43"

# --- de-de: same treatment, its own handback time.
$deMdLink = $null
$deXlfLink = $null
foreach ($h in $wsDe.Hyperlinks) {
    $addr = $h.Range.Address()
    if ($addr -eq '$A$2') { $deMdLink = $h.Address }
    if ($addr -eq '$C$2') { $deXlfLink = $h.Address }
}

$wsDe.Hyperlinks.Add($wsDe.Range("E2"), $deMdLink, $null, $null, "1ddd3d0e-7b9e-45da-b7a2-ededd37b75c7.md")
$wsDe.Hyperlinks.Add($wsDe.Range("F2"), $deXlfLink, $null, $null, "1ddd3d0e-7b9e-45da-b7a2-ededd37b75c7.07b5d910097a9d19b5106c498aa4f82131c0eb33.de-de.xlf")
$wsDe.Range("G2").Value = "2016-01-08 09:21:00"

Write-Host "Handback report generated."
